$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 through 9 (old extra location entries), keep rows 1-3 structure
$ws.Range("A4:E9").EntireRow.Delete() | Out-Null

# Row 2: update to Ringwood / Block 7 Dumplings entry (the "old" exposure period row)
$ws.Range("A2").Value = "Ringwood"
$ws.Range("B2").Value = "Block 7 Dumplings, 171 - 175 Maroondah Highway"
$ws.Range("C2").Value = "29/12/20 8:10am - 8:45am"
$ws.Range("D2").Value = "Case attended store"
$ws.Range("E2").Value = "old"

# Row 3: new corrected exposure period row
$ws.Range("A3").Value = "Ringwood"
$ws.Range("B3").Value = "Block 7 Dumplings, 171 - 175 Maroondah Highway"
$ws.Range("C3").Value = "29/12/20 8:10pm - 8:45pm"
$ws.Range("D3").Value = "Case attended store"
$ws.Range("E3").Value = "new"

# Adjust column widths to fit new (shorter) content
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null

# Select D3, matching the saved selection in the file
$ws.Range("D3").Select() | Out-Null
